$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text values that Excel would otherwise auto-convert to a Number
# (losing the original text/inlineStr cell type) are written via a
# NumberFormat="@" / Style="Normal" bracket so the stored type stays Text
# while the cell keeps its original (default) style index.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '25.971.08'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.638.97'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.60%  '
Set-TextValue $ws.Range('D5') '215.16'
$ws.Range('E5').Value = '  -0.19%  '
Set-TextValue $ws.Range('D6') '0.5135'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('E7').Value = '  -0.57%  '
Set-TextValue $ws.Range('D8') '0.2580'
$ws.Range('E8').Value = '  +0.15%  '
Set-TextValue $ws.Range('D9') '0.06356'
$ws.Range('E9').Value = '  -0.99%  '
Set-TextValue $ws.Range('D10') '19.78'
$ws.Range('E10').Value = '  +0.44%  '
Set-TextValue $ws.Range('D11') '0.07753'
$ws.Range('E11').Value = '  -0.26%  '
Set-TextValue $ws.Range('D12') '4.279'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '1.636.22'
$ws.Range('E13').Value = '  -0.61%  '
Set-TextValue $ws.Range('D14') '0.5463'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '0.0₅7758'
$ws.Range('E15').Value = '  -1.83%  '
Set-TextValue $ws.Range('D16') '64.38'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').Value = '25.998.47'
$ws.Range('E17').Value = '  -0.01%  '
Set-TextValue $ws.Range('D18') '1.003'
$ws.Range('E18').Value = '  -0.36%  '
Set-TextValue $ws.Range('D19') '197.82'
$ws.Range('E19').Value = '  +0.29%  '
Set-TextValue $ws.Range('D20') '4.452'
$ws.Range('E20').Value = '  +0.51%  '
Set-TextValue $ws.Range('D21') '9.950'
$ws.Range('E21').Value = '  -0.91%  '
Set-TextValue $ws.Range('D22') '6.089'
$ws.Range('E22').Value = '  +0.30%  '
Set-TextValue $ws.Range('D23') '1.002'
$ws.Range('E23').Value = '  -0.67%  '
Set-TextValue $ws.Range('D24') '1.913'
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('E26').Value = '  +7.93%  '
Set-TextValue $ws.Range('D27') '6.854'
$ws.Range('E27').Value = '  -0.64%  '
Set-TextValue $ws.Range('D28') '15.59'
$ws.Range('E28').Value = '  -1.12%  '
Set-TextValue $ws.Range('D29') '1.241'
$ws.Range('E29').Value = '  +0.02%  '
Set-TextValue $ws.Range('D30') '0.04846'
$ws.Range('E30').Value = '  -3.52%  '
Set-TextValue $ws.Range('D31') '3.290'
$ws.Range('E31').Value = '  +0.50%  '
Set-TextValue $ws.Range('D32') '3.222'
$ws.Range('E32').Value = '  +0.67%  '
Set-TextValue $ws.Range('D33') '1.538'
$ws.Range('E33').Value = '  -0.46%  '
Set-TextValue $ws.Range('D34') '2.375'
$ws.Range('E34').Value = '  +0.18%  '
Set-TextValue $ws.Range('D35') '0.9146'
$ws.Range('E35').Value = '  +2.20%  '
$ws.Range('D36').Value = '1.148.96'
$ws.Range('E36').Value = '  +1.42%  '
Set-TextValue $ws.Range('D37') '0.5568'
$ws.Range('E37').Value = '  +0.21%  '
Set-TextValue $ws.Range('D38') '2.566'
$ws.Range('E38').Value = '  -1.19%  '
Set-TextValue $ws.Range('D39') '0.01571'
$ws.Range('E39').Value = '  +0.49%  '
Set-TextValue $ws.Range('D40') '1.000'
$ws.Range('E40').Value = '  -0.78%  '
Set-TextValue $ws.Range('D41') '2.519'
$ws.Range('E41').Value = '  -2.12%  '
Set-TextValue $ws.Range('D42') '5.574'
$ws.Range('E42').Value = '  -1.73%  '
Set-TextValue $ws.Range('D43') '0.8072'
$ws.Range('E43').Value = '  -1.16%  '
Set-TextValue $ws.Range('D44') '99.47'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('E45').Value = '  -3.21%  '
$ws.Range('D46').Value = '1.780.88'
$ws.Range('E46').Value = '  -0.26%  '
Set-TextValue $ws.Range('D47') '0.4535'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D48') '1.006'
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D49') '55.18'
$ws.Range('E49').Value = '  -0.24%  '
Set-TextValue $ws.Range('D50') '0.05212'
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '7.493'
$ws.Range('E51').Value = '  +1.28%  '
